$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row (row 2) describing collection item MCH247, matching the
# header row: identifier | alternativeIdentifiers | title | date_s |
# levelOfDescription | extentAndMedium | notes | file_path
$ws.Range("A2").Value = "MCH247"
$ws.Range("C2").Value = "PERMIT IN TERMS OF SECTION 21 (GROUP AREAS ACT), MESSAGE FROM BRIGADIER CIA SWART"
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: CABINET 1C | GRAP COUNT NUMER: NONE"

# D2 and H2 stay blank but still belong to the formatted data row.
$ws.Range("D2").Value = ""
$ws.Range("H2").Value = ""

# Row 2 uses a plain (non-bold, non-shaded) 10pt Calibri style, distinct
# from the Arial default and from the bold/shaded header row.
$ws.Range("A2").Font.Name = "Calibri"
$ws.Range("C2:H2").Font.Name = "Calibri"
$ws.Range("A2").Font.Size = 10
$ws.Range("C2:H2").Font.Size = 10

# Select the newly-populated row and re-assert the frozen header pane,
# matching the saved view state of the source workbook.
$ws.Range("A2:H2").Select()
$ws.Application.ActiveWindow.FreezePanes = $true
